$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Locate the two quoted single-letter symbols inside the sentence that
# explains which player gets which symbol, e.g.:
#   ...p simbolo "X", este deve ... sera atribuido como "O" que, ja...
# After the edit they must be swapped:
#   ...p simbolo "O", este deve ... sera atribuido como "X" que, ja...
# ----------------------------------------------------------------------

# --- First symbol (currently X, becomes O) ---------------------------
$r = $d.Content
$found1 = $r.Find.Execute("símbolo “X", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$p2 = $r.End            # just after the X, before the closing curly quote
$p1 = $p2 - 1            # just before the X, after the opening curly quote

# --- Second symbol (currently O, becomes X) ---------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("como “O", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$p4 = $r2.End            # just after the O, before the closing curly quote
$p3 = $p4 - 1            # just before the O, after the opening curly quote

# ----------------------------------------------------------------------
# Force the run to split at each of these four positions by briefly
# dropping (and immediately removing) a bookmark there. The bookmark
# itself disappears, but the run boundary it creates is retained.
# ----------------------------------------------------------------------
$d.Bookmarks.Add("tmpSplit1", $d.Range($p1, $p1))
$d.Bookmarks.Add("tmpSplit2", $d.Range($p2, $p2))
$d.Bookmarks.Add("tmpSplit3", $d.Range($p3, $p3))
$d.Bookmarks.Add("tmpSplit4", $d.Range($p4, $p4))
$d.Bookmarks("tmpSplit1").Delete()
$d.Bookmarks("tmpSplit2").Delete()
$d.Bookmarks("tmpSplit3").Delete()
$d.Bookmarks("tmpSplit4").Delete()

# ----------------------------------------------------------------------
# Now swap the letters. Because both replacements are a single character
# for a single character, none of the positions computed above shift.
# ----------------------------------------------------------------------
$d.Range($p1, $p2).Text = "O"
$d.Range($p3, $p4).Text = "X"

# ----------------------------------------------------------------------
# Move the "_GoBack" bookmark: delete it from its old location (end of
# the "$(símbolo): X: Y" paragraph) and re-create it, collapsed, right
# after the newly placed X (i.e. at position p4).
# ----------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($p4, $p4))
